$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 6.527126666666668
$ws.Range("H2").Value = 19.58138
$ws.Range("I2").Value = 0.2130391554800433
$ws.Range("J2").Value = 0.2130391554800433
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 85.89497033333333
$ws.Range("N2").Value = 257.684911
$ws.Range("O2").Value = 0.7848889718219874
$ws.Range("P2").Value = 0.7848889718219874
$ws.Range("Q2").Value = 560.6473513952423
$ws.Range("R2").Value = 5045.82616255718
$ws.Range("S2").Value = 0.1672120837025557
$ws.Range("T2").Value = 0.1672120837025557
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 6.527126666666668
$ws.Range("H3").Value = 19.58138
$ws.Range("I3").Value = 0.2130391554800433
$ws.Range("J3").Value = 0.2130391554800433
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.726772333333333
$ws.Range("N3").Value = 5.180317
$ws.Range("O3").Value = 0.0157788582500353
$ws.Range("P3").Value = 0.0157788582500353
$ws.Range("Q3").Value = 11.27086174416222
$ws.Range("R3").Value = 101.43775569746
$ws.Range("S3").Value = 0.003361514636026834
$ws.Range("T3").Value = 0.003361514636026834
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 6.527126666666668
$ws.Range("H4").Value = 19.58138
$ws.Range("I4").Value = 0.2130391554800433
$ws.Range("J4").Value = 0.2130391554800433
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 18.21376166666667
$ws.Range("N4").Value = 54.641285
$ws.Range("O4").Value = 0.1664332685846793
$ws.Range("P4").Value = 0.1664332685846793
$ws.Range("Q4").Value = 118.8835294748111
$ws.Range("R4").Value = 1069.9517652733
$ws.Range("S4").Value = 0.0354568029830633
$ws.Range("T4").Value = 0.03545680298306331
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 6.527126666666668
$ws.Range("H5").Value = 19.58138
$ws.Range("I5").Value = 0.2130391554800433
$ws.Range("J5").Value = 0.2130391554800433
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 3.600318333333334
$ws.Range("N5").Value = 10.800955
$ws.Range("O5").Value = 0.03289890134329811
$ws.Range("P5").Value = 0.03289890134329811
$ws.Range("Q5").Value = 23.49973380198889
$ws.Range("R5").Value = 211.4976042179
$ws.Range("S5").Value = 0.007008754158397492
$ws.Range("T5").Value = 0.007008754158397492
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 10.77811266666667
$ws.Range("H6").Value = 32.334338
$ws.Range("I6").Value = 0.351787262211666
$ws.Range("J6").Value = 0.351787262211666
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 85.89497033333333
$ws.Range("N6").Value = 257.684911
$ws.Range("O6").Value = 0.7848889718219874
$ws.Range("P6").Value = 0.7848889718219874
$ws.Range("Q6").Value = 925.7856677526576
$ws.Range("R6").Value = 8332.071009773919
$ws.Range("S6").Value = 0.2761139425373864
$ws.Range("T6").Value = 0.2761139425373864
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 10.77811266666667
$ws.Range("H7").Value = 32.334338
$ws.Range("I7").Value = 0.351787262211666
$ws.Range("J7").Value = 0.351787262211666
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.726772333333333
$ws.Range("N7").Value = 5.180317
$ws.Range("O7").Value = 0.0157788582500353
$ws.Range("P7").Value = 0.0157788582500353
$ws.Range("Q7").Value = 18.61134675834955
$ws.Range("R7").Value = 167.502120825146
$ws.Range("S7").Value = 0.005550801344605877
$ws.Range("T7").Value = 0.005550801344605877
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 10.77811266666667
$ws.Range("H8").Value = 32.334338
$ws.Range("I8").Value = 0.351787262211666
$ws.Range("J8").Value = 0.351787262211666
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 18.21376166666667
$ws.Range("N8").Value = 54.641285
$ws.Range("O8").Value = 0.1664332685846793
$ws.Range("P8").Value = 0.1664332685846793
$ws.Range("Q8").Value = 196.3099753271478
$ws.Range("R8").Value = 1766.78977794433
$ws.Range("S8").Value = 0.0585491038963432
$ws.Range("T8").Value = 0.05854910389634321
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 10.77811266666667
$ws.Range("H9").Value = 32.334338
$ws.Range("I9").Value = 0.351787262211666
$ws.Range("J9").Value = 0.351787262211666
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.600318333333334
$ws.Range("N9").Value = 10.800955
$ws.Range("O9").Value = 0.03289890134329811
$ws.Range("P9").Value = 0.03289890134329811
$ws.Range("Q9").Value = 38.80463663253222
$ws.Range("R9").Value = 349.24172969279
$ws.Range("S9").Value = 0.01157341443333054
$ws.Range("T9").Value = 0.01157341443333054
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 6.559571333333333
$ws.Range("H10").Value = 19.678714
$ws.Range("I10").Value = 0.2140981182885632
$ws.Range("J10").Value = 0.2140981182885631
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 85.89497033333333
$ws.Range("N10").Value = 257.684911
$ws.Range("O10").Value = 0.7848889718219874
$ws.Range("P10").Value = 0.7848889718219874
$ws.Range("Q10").Value = 563.4341850760504
$ws.Range("R10").Value = 5070.907665684454
$ws.Range("S10").Value = 0.1680432519325326
$ws.Range("T10").Value = 0.1680432519325326
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 6.559571333333333
$ws.Range("H11").Value = 19.678714
$ws.Range("I11").Value = 0.2140981182885632
$ws.Range("J11").Value = 0.2140981182885631
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 1.726772333333333
$ws.Range("N11").Value = 5.180317
$ws.Range("O11").Value = 0.0157788582500353
$ws.Range("P11").Value = 0.0157788582500353
$ws.Range("Q11").Value = 11.32688629692644
$ws.Range("R11").Value = 101.941976672338
$ws.Range("S11").Value = 0.003378223860074528
$ws.Range("T11").Value = 0.003378223860074528
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 6.559571333333333
$ws.Range("H12").Value = 19.678714
$ws.Range("I12").Value = 0.2140981182885632
$ws.Range("J12").Value = 0.2140981182885631
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 18.21376166666667
$ws.Range("N12").Value = 54.641285
$ws.Range("O12").Value = 0.1664332685846793
$ws.Range("P12").Value = 0.1664332685846793
$ws.Range("Q12").Value = 119.4744689008322
$ws.Range("R12").Value = 1075.27022010749
$ws.Range("S12").Value = 0.03563304962459487
$ws.Range("T12").Value = 0.03563304962459487
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 6.559571333333333
$ws.Range("H13").Value = 19.678714
$ws.Range("I13").Value = 0.2140981182885632
$ws.Range("J13").Value = 0.2140981182885631
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 3.600318333333334
$ws.Range("N13").Value = 10.800955
$ws.Range("O13").Value = 0.03289890134329811
$ws.Range("P13").Value = 0.03289890134329811
$ws.Range("Q13").Value = 23.61654493020778
$ws.Range("R13").Value = 212.54890437187
$ws.Range("S13").Value = 0.007043592871361208
$ws.Range("T13").Value = 0.007043592871361207
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 6.773344333333334
$ws.Range("H14").Value = 20.320033
$ws.Range("I14").Value = 0.2210754640197275
$ws.Range("J14").Value = 0.2210754640197275
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 85.89497033333333
$ws.Range("N14").Value = 257.684911
$ws.Range("O14").Value = 0.7848889718219874
$ws.Range("P14").Value = 0.7848889718219874
$ws.Range("Q14").Value = 581.7962105691182
$ws.Range("R14").Value = 5236.165895122063
$ws.Range("S14").Value = 0.1735196936495127
$ws.Range("T14").Value = 0.1735196936495126
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 6.773344333333334
$ws.Range("H15").Value = 20.320033
$ws.Range("I15").Value = 0.2210754640197275
$ws.Range("J15").Value = 0.2210754640197275
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 1.726772333333333
$ws.Range("N15").Value = 5.180317
$ws.Range("O15").Value = 0.0157788582500353
$ws.Range("P15").Value = 0.0157788582500353
$ws.Range("Q15").Value = 11.69602359894011
$ws.Range("R15").Value = 105.264212390461
$ws.Range("S15").Value = 0.003488318409328059
$ws.Range("T15").Value = 0.003488318409328058
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 6.773344333333334
$ws.Range("H16").Value = 20.320033
$ws.Range("I16").Value = 0.2210754640197275
$ws.Range("J16").Value = 0.2210754640197275
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 18.21376166666667
$ws.Range("N16").Value = 54.641285
$ws.Range("O16").Value = 0.1664332685846793
$ws.Range("P16").Value = 0.1664332685846793
$ws.Range("Q16").Value = 123.3680793736006
$ws.Range("R16").Value = 1110.312714362405
$ws.Range("S16").Value = 0.0367943120806779
$ws.Range("T16").Value = 0.0367943120806779
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 6.773344333333334
$ws.Range("H17").Value = 20.320033
$ws.Range("I17").Value = 0.2210754640197275
$ws.Range("J17").Value = 0.2210754640197275
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 3.600318333333334
$ws.Range("N17").Value = 10.800955
$ws.Range("O17").Value = 0.03289890134329811
$ws.Range("P17").Value = 0.03289890134329811
$ws.Range("Q17").Value = 24.38619578127945
$ws.Range("R17").Value = 219.475762031515
$ws.Range("S17").Value = 0.007273139880208866
$ws.Range("T17").Value = 0.007273139880208865
